$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update weekly price data for Hortaliza - Alcachofa (Mapocho Venta Directa de Santiago)
$ws.Cells.Item(2, 4).Value = 44435
$ws.Cells.Item(2, 10).Value = 25
$ws.Cells.Item(2, 11).Value = 14000
$ws.Cells.Item(2, 12).Value = 14000
$ws.Cells.Item(2, 13).Value = 14000
$ws.Cells.Item(2, 16).Value = 467
$ws.Cells.Item(3, 4).Value = 44435
$ws.Cells.Item(3, 10).Value = 25
$ws.Cells.Item(3, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(4, 4).Value = 44432
$ws.Cells.Item(4, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(5, 4).Value = 44467
$ws.Cells.Item(5, 10).Value = 35
$ws.Cells.Item(6, 4).Value = 44418
$ws.Cells.Item(6, 10).Value = 30
$ws.Cells.Item(6, 11).Value = 15000
$ws.Cells.Item(6, 12).Value = 15000
$ws.Cells.Item(6, 13).Value = 15000
$ws.Cells.Item(6, 16).Value = 500
$ws.Cells.Item(7, 4).Value = 44474
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(7, 16).Value = 333
$ws.Cells.Item(7, 17).Value = 30
$ws.Cells.Item(8, 4).Value = 44446
$ws.Cells.Item(8, 10).Value = 25
$ws.Cells.Item(8, 11).Value = 14000
$ws.Cells.Item(8, 12).Value = 14000
$ws.Cells.Item(8, 13).Value = 14000
$ws.Cells.Item(8, 16).Value = 467
$ws.Cells.Item(9, 4).Value = 44460
$ws.Cells.Item(9, 10).Value = 45
$ws.Cells.Item(9, 11).Value = 13000
$ws.Cells.Item(9, 12).Value = 13000
$ws.Cells.Item(9, 13).Value = 13000
$ws.Cells.Item(9, 16).Value = 433
$ws.Cells.Item(10, 4).Value = 44376
$ws.Cells.Item(10, 11).Value = 18000
$ws.Cells.Item(10, 12).Value = 18000
$ws.Cells.Item(10, 13).Value = 18000
$ws.Cells.Item(10, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 16).Value = 600
$ws.Cells.Item(11, 4).Value = 44453
$ws.Cells.Item(11, 10).Value = 50
$ws.Cells.Item(11, 11).Value = 12000
$ws.Cells.Item(11, 12).Value = 12000
$ws.Cells.Item(11, 13).Value = 12000
$ws.Cells.Item(11, 16).Value = 400
$ws.Cells.Item(12, 4).Value = 44421
$ws.Cells.Item(12, 11).Value = 15000
$ws.Cells.Item(12, 12).Value = 16000
$ws.Cells.Item(12, 13).Value = 15400
$ws.Cells.Item(12, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(12, 16).Value = 513
$ws.Cells.Item(13, 4).Value = 44841
$ws.Cells.Item(13, 10).Value = 45
$ws.Cells.Item(13, 11).Value = 12000
$ws.Cells.Item(13, 12).Value = 12000
$ws.Cells.Item(13, 13).Value = 12000
$ws.Cells.Item(13, 16).Value = 400
$ws.Cells.Item(14, 4).Value = 44841
$ws.Cells.Item(14, 9).Value = "Segunda"
$ws.Cells.Item(14, 11).Value = 10000
$ws.Cells.Item(14, 12).Value = 10000
$ws.Cells.Item(14, 13).Value = 10000
$ws.Cells.Item(14, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(14, 16).Value = 250
$ws.Cells.Item(14, 17).Value = 40
$ws.Cells.Item(15, 4).Value = 44449
$ws.Cells.Item(15, 10).Value = 45
$ws.Cells.Item(15, 11).Value = 12000
$ws.Cells.Item(15, 12).Value = 12000
$ws.Cells.Item(15, 13).Value = 12000
$ws.Cells.Item(15, 16).Value = 400
$ws.Cells.Item(16, 4).Value = 44425
$ws.Cells.Item(16, 11).Value = 14000
$ws.Cells.Item(16, 12).Value = 14000
$ws.Cells.Item(16, 13).Value = 14000
$ws.Cells.Item(16, 16).Value = 467
